$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F2").Value = 14088
$ws1.Range("F6").Value = 661
$ws1.Range("F11").Value = 2177
$ws1.Range("F12").Value = 168
$ws1.Range("F13").Value = 120
$ws1.Range("F14").Value = 96
$ws1.Range("F15").Value = 210
$ws1.Range("F17").Value = 587
$ws1.Range("F19").Value = 499
$ws1.Range("F21").Value = 30
$ws1.Range("F22").Value = 303
$ws1.Range("F23").Value = 883
$ws1.Range("F24").Value = 141
$ws1.Range("F25").Value = 70
$ws1.Range("F26").Value = 24
$ws1.Range("F29").Value = 70
$ws1.Range("F30").Value = 32
$ws2.Range("F7").Value = 182
$ws2.Range("F8").Value = 1670
$ws2.Range("F15").Value = 1812
$ws3.Range("F3").Value = 134
$ws4.Range("F3").Value = 14088
$ws4.Range("F7").Value = 661
$ws4.Range("F14").Value = 2177
$ws4.Range("F15").Value = 134
$ws4.Range("F16").Value = 168
$ws4.Range("F17").Value = 168
$ws4.Range("F18").Value = 120
$ws4.Range("F19").Value = 96
$ws4.Range("F20").Value = 210
$ws4.Range("F26").Value = 587
$ws4.Range("F28").Value = 499
$ws4.Range("F30").Value = 30
$ws4.Range("F31").Value = 303
$ws4.Range("F32").Value = 883
$ws4.Range("F33").Value = 182
$ws4.Range("F34").Value = 1670
$ws4.Range("F39").Value = 141
$ws4.Range("F40").Value = 70
$ws4.Range("F41").Value = 24
$ws4.Range("F46").Value = 70
$ws4.Range("F47").Value = 32
$ws4.Range("F48").Value = 1812
